# Update "Summary of model after N samples" sections:
#  - rename/renumber the first 8 sections and refresh their stats
#  - drop the trailing sections (37, 38, 39, 40, 42, 43, 44) entirely

$d = $word.ActiveDocument

$updates = @(
    @{ oldN = 28; newN = 209; oldAcc = "95.92%"; newAcc = "93.17%"; oldLoss = "4.08%"; newLoss = "6.83%" },
    @{ oldN = 30; newN = 210; oldAcc = "92.65%"; newAcc = "93.18%"; oldLoss = "7.35%"; newLoss = "6.82%" },
    @{ oldN = 31; newN = 211; oldAcc = "92.63%"; newAcc = "93.19%"; oldLoss = "7.37%"; newLoss = "6.81%" },
    @{ oldN = 32; newN = 294; oldAcc = "92.81%"; newAcc = "93.89%"; oldLoss = "7.19%"; newLoss = "6.11%" },
    @{ oldN = 33; newN = 295; oldAcc = "92.98%"; newAcc = "93.90%"; oldLoss = "7.02%"; newLoss = "6.10%" },
    @{ oldN = 34; newN = 296; oldAcc = "93.14%"; newAcc = "93.91%"; oldLoss = "6.86%"; newLoss = "6.09%" },
    @{ oldN = 35; newN = 297; oldAcc = "93.30%"; newAcc = "93.93%"; oldLoss = "6.70%"; newLoss = "6.07%" },
    @{ oldN = 36; newN = 298; oldAcc = "93.30%"; newAcc = "93.94%"; oldLoss = "6.70%"; newLoss = "6.06%" }
)

# --- 1. Refresh the stats inside the first eight tables (tables keep their
#        document order, so table i is the one belonging to $updates[i-1]) --

for ($i = 0; $i -lt $updates.Count; $i++) {
    $u = $updates[$i]
    $t = $d.Tables($i + 1)
    $rng = $t.Range
    $rng.Find.Execute([string]$u.oldN, $true, $false, $false, $false, $false,
                       $true, 1, $false, [string]$u.newN, 2) | Out-Null
    $rng = $t.Range
    $rng.Find.Execute($u.oldAcc, $true, $false, $false, $false, $false,
                       $true, 1, $false, $u.newAcc, 2) | Out-Null
    $rng = $t.Range
    $rng.Find.Execute($u.oldLoss, $true, $false, $false, $false, $false,
                       $true, 1, $false, $u.newLoss, 2) | Out-Null
}

# --- 2. Rename the Heading1 paragraphs that introduce each section --------

foreach ($u in $updates) {
    $d.Content.Find.Execute("Summary of model after $($u.oldN) samples", $true, $false, $false, $false, $false,
                             $true, 1, $false, "Summary of model after $($u.newN) samples", 2) | Out-Null
}

# --- 3. Remove the sections that no longer exist ---------------------------

$toRemove = @(44, 43, 42, 40, 39, 38, 37)
foreach ($n in $toRemove) {
    $heading = "Summary of model after $n samples"
    $foundAt = -1
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.TrimEnd() -eq $heading) {
            $foundAt = $i
            break
        }
    }
    if ($foundAt -ge 1) {
        $headingPara = $d.Paragraphs($foundAt)
        $headingEnd = $headingPara.Range.End
        for ($j = 1; $j -le $d.Tables.Count; $j++) {
            $t = $d.Tables($j)
            if ($t.Range.Start -eq $headingEnd) {
                $t.Delete()
                break
            }
        }
        $headingPara.Range.Delete()
    }
}
